$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Data4" column header
$ws.Range("E1").Value = "Data4"

# Replace the old "tc001_Login" row (row 2) with the new "tc001_addBook" book row
$ws.Range("B2").Value = "Jack"
$ws.Range("A2").Value = "tc001_addBook"
$ws.Range("E2").Value = "J.K Rowling"
$ws.Range("C2").Value = "ISBN:4679-"
$ws.Range("D2").Value = 5435

# New column widths for the newly-relevant columns (closest the host's
# pixel-quantized ColumnWidth setter can land to the authored 12.109375 /
# 11.88671875 / 15.33203125 stored widths)
$ws.Columns.Item(3).ColumnWidth = 11.25
$ws.Columns.Item(4).ColumnWidth = 10.916666666666666
$ws.Columns.Item(5).ColumnWidth = 14.416666666666666

# Move the selection like the author's session ended up
$ws.Range("E12").Select()
